$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.456.01'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').Value = '1.827.35'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = "'315.77"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').Value = "'0.5167"
$ws.Range('E7').Value = '  -4.16%  '
$ws.Range('D8').Value = "'0.3936"
$ws.Range('E8').Value = '  +2.92%  '
$ws.Range('D9').Value = "'0.07721"
$ws.Range('E9').Value = '  +3.99%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  +2.54%  '
$ws.Range('D12').Value = "'21.08"
$ws.Range('E12').Value = '  +3.95%  '
$ws.Range('D13').Value = "'6.283"
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '1.822.86'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').Value = "'93.28"
$ws.Range('E17').Value = '  +5.47%  '
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').Value = "'0.06632"
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').Value = "'17.71"
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').Value = "'6.066"
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('D23').Value = '28.461.53'
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('D24').Value = "'11.14"
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').Value = "'2.259"
$ws.Range('E25').Value = '  +8.16%  '
$ws.Range('D26').Value = "'2.452"
$ws.Range('E26').Value = '  +4.97%  '
$ws.Range('D27').Value = "'157.32"
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('E28').Value = '  +1.85%  '
$ws.Range('D29').Value = '2.036.00'
$ws.Range('E29').Value = '  +1.65%  '
$ws.Range('D30').Value = "'124.97"
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').Value = "'1.129"
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('D32').Value = "'0.1100"
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = "'5.658"
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('D34').Value = "'3.674"
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').Value = "'0.07184"
$ws.Range('E35').Value = '  +3.39%  '
$ws.Range('D36').Value = "'0.2235"
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').Value = "'8.967"
$ws.Range('E37').Value = '  +5.61%  '
$ws.Range('D38').Value = "'0.02323"
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('D39').Value = "'5.153"
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').Value = "'0.6246"
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').Value = "'1.193"
$ws.Range('E42').Value = '  +2.56%  '
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').Value = "'1.394"
$ws.Range('E44').Value = '  -1.54%  '
$ws.Range('D45').Value = "'13.47"
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('D46').Value = "'0.5906"
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('D47').Value = "'3.705"
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').Value = "'124.38"
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = "'1.980"
$ws.Range('E49').Value = '  +3.78%  '
$ws.Range('D50').Value = "'1.184"
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').Value = "'0.06925"
$ws.Range('E51').Value = '  +2.01%  '
